$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove all existing hyperlinks on the sheet so we can rebuild them cleanly
$ws.Range("F2").Hyperlinks.Delete()

# Row 2
$ws.Range("A2").Value = '2026-01-17 01:23:54'
$ws.Range("B2").Value = '製造業向け図面自動生成システムの開発・ツール化を支援してくださるエンジニア募集(AI/バックエンド)'
$ws.Range("C2").Value = 'システム開発'
$ws.Range("D2").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E2").Value = '期限情報なし'
$ws.Range("F2").Value = 'https://www.lancers.jp/work/detail/5460562'
$ws.Range("G2").Value = 435
$ws.Range("H2").Value = '🔥AI,Ai ◆ツール,開発'

# Row 3
$ws.Range("A3").Value = '2026-01-17 01:23:54'
$ws.Range("B3").Value = '【募集】Python / Docker 日次データ スクレイピングシステム構築'
$ws.Range("C3").Value = 'システム開発'
$ws.Range("D3").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E3").Value = '期限情報なし'
$ws.Range("F3").Value = 'https://www.lancers.jp/work/detail/5469627'
$ws.Range("G3").Value = 248
$ws.Range("H3").Value = '🔥Python ◆スクレイピング'

# Row 4
$ws.Range("A4").Value = '2026-01-17 01:23:54'
$ws.Range("B4").Value = '※急募:Next.jsによる業務アプリの開発(+Flutter)'
$ws.Range("C4").Value = 'システム開発'
$ws.Range("D4").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E4").Value = '期限情報なし'
$ws.Range("F4").Value = 'https://www.lancers.jp/work/detail/5473147'
$ws.Range("G4").Value = 225
$ws.Range("H4").Value = '🔥Next.js ◆開発 ◇アプリ'

# Row 5
$ws.Range("A5").Value = '2026-01-17 01:23:54'
$ws.Range("B5").Value = '施設管理・現場業務向け チェックリスト業務の自動化・報告書作成システム開発エンジニア募集'
$ws.Range("C5").Value = 'システム開発'
$ws.Range("D5").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E5").Value = '期限情報なし'
$ws.Range("F5").Value = 'https://www.lancers.jp/work/detail/5460563'
$ws.Range("G5").Value = 220
$ws.Range("H5").Value = '◆開発,システム開発 ◇管理'

# Row 6
$ws.Range("A6").Value = '2026-01-17 01:23:54'
$ws.Range("B6").Value = '※急募:Flutterによる業務アプリの開発(+next.js)'
$ws.Range("C6").Value = 'システム開発'
$ws.Range("D6").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E6").Value = '期限情報なし'
$ws.Range("F6").Value = 'https://www.lancers.jp/work/detail/5473146'
$ws.Range("G6").Value = 218
$ws.Range("H6").Value = '🔥Next.js ◆開発 ◇アプリ'

# Row 7
$ws.Range("A7").Value = '2026-01-17 01:23:54'
$ws.Range("B7").Value = '【募集】RPAツール「RoboTANGO」設定代行の専門家を探しています'
$ws.Range("C7").Value = 'システム開発'
$ws.Range("D7").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E7").Value = '期限情報なし'
$ws.Range("F7").Value = 'https://www.lancers.jp/work/detail/5405023'
$ws.Range("G7").Value = 178
$ws.Range("H7").Value = '★bot ◆ツール'

# Row 8
$ws.Range("A8").Value = '2026-01-17 01:23:54'
$ws.Range("B8").Value = '【急募】Accessでの受発注管理・請求書発行システム開発'
$ws.Range("C8").Value = 'システム開発'
$ws.Range("D8").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E8").Value = '期限情報なし'
$ws.Range("F8").Value = 'https://www.lancers.jp/work/detail/5473234'
$ws.Range("G8").Value = 148
$ws.Range("H8").Value = '◆開発,システム開発 ◇管理'

# Row 9
$ws.Range("A9").Value = '2026-01-17 01:23:54'
$ws.Range("B9").Value = '【バイナリ解析 / 逆コンパイル】EPCデータ解析ツール開発'
$ws.Range("C9").Value = 'システム開発'
$ws.Range("D9").Value = '1,000,000 円 ~ 3,000,000 円 / 固定'
$ws.Range("E9").Value = '期限情報なし'
$ws.Range("F9").Value = 'https://www.lancers.jp/work/detail/5473181'
$ws.Range("G9").Value = 135
$ws.Range("H9").Value = '◆ツール,開発'

# Row 10
$ws.Range("A10").Value = '2026-01-17 01:23:54'
$ws.Range("B10").Value = '【Windows/Wacom】署名画像から筆順解析図を作成する業務用アプリ開発'
$ws.Range("C10").Value = 'システム開発'
$ws.Range("D10").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E10").Value = '期限情報なし'
$ws.Range("F10").Value = 'https://www.lancers.jp/work/detail/5472804'
$ws.Range("G10").Value = 100
$ws.Range("H10").Value = '◆開発 ◇アプリ'

# Row 11
$ws.Range("A11").Value = '2026-01-17 01:23:54'
$ws.Range("B11").Value = '署名画像から筆順を可視化するアプリ開発者募集'
$ws.Range("C11").Value = 'システム開発'
$ws.Range("D11").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E11").Value = '期限情報なし'
$ws.Range("F11").Value = 'https://www.lancers.jp/work/detail/5472080'
$ws.Range("G11").Value = 100
$ws.Range("H11").Value = '◆開発 ◇アプリ'

# Row 12
$ws.Range("A12").Value = '2026-01-17 01:23:54'
$ws.Range("B12").Value = 'スマホカラオケ予約Webアプリ開発のフリーランス募集(使用するのは個人の集まりで趣味で使う程度です)'
$ws.Range("C12").Value = 'システム開発'
$ws.Range("D12").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E12").Value = '期限情報なし'
$ws.Range("F12").Value = 'https://www.lancers.jp/work/detail/5472431'
$ws.Range("G12").Value = 88
$ws.Range("H12").Value = '◆開発 ◇アプリ'

# Row 13
$ws.Range("A13").Value = '2026-01-17 01:23:54'
$ws.Range("B13").Value = '【自動化】申込書AからBへの転写をエクセルで実現したい'
$ws.Range("C13").Value = 'システム開発'
$ws.Range("D13").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E13").Value = '期限情報なし'
$ws.Range("F13").Value = 'https://www.lancers.jp/work/detail/5473042'
$ws.Range("G13").Value = 83
$ws.Range("H13").Value = '◆自動化'

# Row 14
$ws.Range("A14").Value = '2026-01-17 01:23:54'
$ws.Range("B14").Value = '初回 Webサーバ管理エンジニア'
$ws.Range("C14").Value = 'システム開発'
$ws.Range("D14").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E14").Value = '期限情報なし'
$ws.Range("F14").Value = 'https://www.lancers.jp/work/detail/5472544'
$ws.Range("G14").Value = 45
$ws.Range("H14").Value = '◇管理'

# Row 15
$ws.Range("A15").Value = '2026-01-17 01:23:54'
$ws.Range("B15").Value = '【急募】Flutterflowの扱えるノーコードエンジニアを探しています!'
$ws.Range("C15").Value = 'システム開発'
$ws.Range("D15").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E15").Value = '期限情報なし'
$ws.Range("F15").Value = 'https://www.lancers.jp/work/detail/5472976'
$ws.Range("G15").Value = 25
$ws.Range("H15").ClearContents()

# Row 16
$ws.Range("A16").Value = '2026-01-17 01:23:54'
$ws.Range("B16").Value = 'm.2 SSD基板の設計'
$ws.Range("C16").Value = 'システム開発'
$ws.Range("D16").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E16").Value = '期限情報なし'
$ws.Range("F16").Value = 'https://www.lancers.jp/work/detail/5472120'
$ws.Range("G16").Value = 25
$ws.Range("H16").ClearContents()

# Row 17
$ws.Range("A17").Value = '2026-01-17 01:23:54'
$ws.Range("B17").Value = '《長期レギュラー》公的機関Web運用の要となる、ディレクター募集'
$ws.Range("C17").Value = 'システム開発'
$ws.Range("D17").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E17").Value = '期限情報なし'
$ws.Range("F17").Value = 'https://www.lancers.jp/work/detail/5472958'
$ws.Range("G17").Value = 18
$ws.Range("H17").ClearContents()

# Re-create hyperlinks for F2:F17 and apply Hyperlink style
$ws.Hyperlinks.Add($ws.Range("F2"), 'https://www.lancers.jp/work/detail/5460562')
$ws.Hyperlinks.Add($ws.Range("F3"), 'https://www.lancers.jp/work/detail/5469627')
$ws.Hyperlinks.Add($ws.Range("F4"), 'https://www.lancers.jp/work/detail/5473147')
$ws.Hyperlinks.Add($ws.Range("F5"), 'https://www.lancers.jp/work/detail/5460563')
$ws.Hyperlinks.Add($ws.Range("F6"), 'https://www.lancers.jp/work/detail/5473146')
$ws.Hyperlinks.Add($ws.Range("F7"), 'https://www.lancers.jp/work/detail/5405023')
$ws.Hyperlinks.Add($ws.Range("F8"), 'https://www.lancers.jp/work/detail/5473234')
$ws.Hyperlinks.Add($ws.Range("F9"), 'https://www.lancers.jp/work/detail/5473181')
$ws.Hyperlinks.Add($ws.Range("F10"), 'https://www.lancers.jp/work/detail/5472804')
$ws.Hyperlinks.Add($ws.Range("F11"), 'https://www.lancers.jp/work/detail/5472080')
$ws.Hyperlinks.Add($ws.Range("F12"), 'https://www.lancers.jp/work/detail/5472431')
$ws.Hyperlinks.Add($ws.Range("F13"), 'https://www.lancers.jp/work/detail/5473042')
$ws.Hyperlinks.Add($ws.Range("F14"), 'https://www.lancers.jp/work/detail/5472544')
$ws.Hyperlinks.Add($ws.Range("F15"), 'https://www.lancers.jp/work/detail/5472976')
$ws.Hyperlinks.Add($ws.Range("F16"), 'https://www.lancers.jp/work/detail/5472120')
$ws.Hyperlinks.Add($ws.Range("F17"), 'https://www.lancers.jp/work/detail/5472958')
$ws.Range("F2:F17").Style = "Hyperlink"

# Column width adjustments (D: 30 -> 32, H: 18 -> 19)
$ws.Columns.Item(4).ColumnWidth = 31.17
$ws.Columns.Item(8).ColumnWidth = 18.17

